# IntroPandas.xlsx — add the "2885. Rename Columns" LeetCode pandas entry
# as a new row at the bottom of the existing Table2 listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row; this keeps the table ref / autoFilter / dimension
# consistent with the rest of the sheet.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()

# Fill in the new row's cells (row 10), matching the column layout:
# Question | Difficulty | Pattern | Notes | Link
$ws.Range("A10").Value = "2885. Rename Columns"
$ws.Range("B10").Value = "Easy"
$ws.Range("C10").Value = "Data Cleaning"
$ws.Range("D10").Value = "Use .rename() method. "
$ws.Range("E10").Value = "https://leetcode.com/problems/rename-columns/solutions/4141051/line-by-line-explanation-easy-solution-beginner-friendly-pandas/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata "

# Match the green "Easy" difficulty highlight used by the other rows.
$ws.Range("B10").Interior.Color = $ws.Range("B9").Interior.Color

# Turn the Link cell into a real hyperlink, styled like the other Link cells.
$ws.Hyperlinks.Add($ws.Range("E10"), "https://leetcode.com/problems/rename-columns/solutions/4141051/line-by-line-explanation-easy-solution-beginner-friendly-pandas/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata ")
$ws.Range("E10").Style = "Hyperlink"

# Update the saved selection/active cell like the source workbook.
$ws.Range("D23").Select()
